$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.379.35'
$ws.Range("E2").Value = '  -0.83%  '

$ws.Range("D3").Value = '1.637.57'
$ws.Range("E3").Value = '  -1.71%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.09'
$ws.Range("E5").Value = '  -1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.530'
$ws.Range("E6").Value = '  +3.68%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.04'

$ws.Range("E9").Value = '  -2.61%  '

$ws.Range("E10").Value = '  -2.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +1.27%  '

$ws.Range("D12").Value = '1.869.22'
$ws.Range("E12").Value = '  -1.69%  '

$ws.Range("D13").Value = '1.623.34'
$ws.Range("E13").Value = '  -2.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("E14").Value = '  -2.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.557'
$ws.Range("E15").Value = '  -0.88%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.25'
$ws.Range("E16").Value = '  -2.91%  '

$ws.Range("D17").Value = '27.350.15'
$ws.Range("E17").Value = '  -0.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.42'
$ws.Range("E18").Value = '  -5.65%  '

$ws.Range("E19").Value = '  -1.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.51'
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  -3.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.42'
$ws.Range("E23").Value = '  +1.48%  '

$ws.Range("E24").Value = '  -0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.06'
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.95'
$ws.Range("E26").Value = '  -3.38%  '

$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  -5.51%  '

$ws.Range("E30").Value = '  -4.46%  '

$ws.Range("E31").Value = '  -3.44%  '

$ws.Range("E32").Value = '  -2.40%  '

$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("D34").Value = '1.405.27'
$ws.Range("E34").Value = '  -4.41%  '

$ws.Range("E35").Value = '  +0.11%  '

$ws.Range("E37").Value = '  -1.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.878'
$ws.Range("E38").Value = '  -5.58%  '

$ws.Range("E39").Value = '  -3.64%  '

$ws.Range("E40").Value = '  +0.95%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("E42").Value = '  -1.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.48'
$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("E44").Value = '  +0.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.790'
$ws.Range("E45").Value = '  +0.48%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.37'
$ws.Range("E46").Value = '  -7.20%  '

$ws.Range("D47").Value = '1.779.06'
$ws.Range("E47").Value = '  -1.63%  '

$ws.Range("E48").Value = '  -4.70%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.22'
$ws.Range("E49").Value = '  -2.33%  '

$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  -2.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0988'
$ws.Range("E51").Value = '  -3.91%  '
